$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 10: "plank challenge" results for week of 2022-01-22 (serial 44583) ---

# Reuse the date-formatted style (s="1") from the cell above instead of creating a
# brand-new number-format style entry: copy A9's format down onto A10.
$null = $ws.Range("A9").Copy()
$null = $ws.Range("A10").PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(10, 1).Value = 44583

$ws.Cells.Item(10, 2).Value = 3.93
$ws.Cells.Item(10, 3).Value = 8.2
$ws.Cells.Item(10, 4).Value = 0

$ws.Cells.Item(10, 5).Value = 3.25
$ws.Cells.Item(10, 6).Value = 1.83
$ws.Cells.Item(10, 7).Value = 3.75

$ws.Cells.Item(10, 8).Value = 3.08
$ws.Cells.Item(10, 9).Value = 1.1
$ws.Cells.Item(10, 10).Value = 0.93

$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = 0

$ws.Cells.Item(10, 14).Value = 0
$ws.Cells.Item(10, 15).Value = 0
$ws.Cells.Item(10, 16).Value = 0

$ws.Cells.Item(10, 17).Value = 0
$ws.Cells.Item(10, 18).Value = 0
$ws.Cells.Item(10, 19).Value = 0

$ws.Cells.Item(10, 20).Value = 2.5
$ws.Cells.Item(10, 21).Value = 2
$ws.Cells.Item(10, 22).Value = 3.1

$ws.Cells.Item(10, 23).Formula = "=12*(B10+C10+D10)"
$ws.Cells.Item(10, 24).Formula = "=12*(E10+F10+G10)"
$ws.Cells.Item(10, 25).Formula = "=12*(H10+I10+J10)"
$ws.Cells.Item(10, 26).Formula = "=12*(K10+L10+M10)"
$ws.Cells.Item(10, 27).Formula = "=12*(N10+O10+P10)"
$ws.Cells.Item(10, 28).Formula = "=12*(Q10+R10+S10)"
$ws.Cells.Item(10, 29).Formula = "=12*(T10+U10+V10)"

# --- Column W ("Sum_1") widened to fit the new figures ---
# (engine's char->pixel grid is coarse; this is the closest reachable width to
# the authored 18.1640625)
$ws.Columns.Item(23).ColumnWidth = 17.33

# --- Selection moved to AC11 after data entry ---
$null = $ws.Range("AC11").Select()
